# Weekly price update: insert a new data row for the latest week at the
# top of the price history (row 51, just below the most recent existing
# row 50) and push the older rows down by one. This mirrors Excel's
# native "Insert" behavior, which also copies the row-above formatting
# (needed so the new date cell in column D keeps the date/time style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 51, shifting rows 51:60 down to 52:61.
$ws.Rows(51).Insert()

# Populate the newly inserted row with the latest week's values.
$ws.Cells.Item(51, 1).Value  = 1
$ws.Cells.Item(51, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(51, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(51, 4).Value  = 44642
$ws.Cells.Item(51, 5).Value  = 15
$ws.Cells.Item(51, 6).Value  = 100112012
$ws.Cells.Item(51, 7).Value  = "Espinaca"
$ws.Cells.Item(51, 8).Value  = "Sin especificar"
$ws.Cells.Item(51, 9).Value  = "Primera"
$ws.Cells.Item(51, 10).Value = 250
$ws.Cells.Item(51, 11).Value = 1500
$ws.Cells.Item(51, 12).Value = 2000
$ws.Cells.Item(51, 13).Value = 1750
$ws.Cells.Item(51, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 583
$ws.Cells.Item(51, 17).Value = 3
$ws.Cells.Item(51, 18).Value = "Hortaliza"
